# Roomnames.pptx edit:
#  1) Update the cached "datetimeFigureOut" footer field from 10/10/2024
#     to 12/10/2024 on every slide layout and on the slide master.
#  2) Add a new "RESTRICTED ACCESS" room-name label textbox to slide 2,
#     matching the style of the existing room-name textboxes there.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: 10/10/2024 -> 12/10/2024
# ---------------------------------------------------------------------
$oldDate = "10/10/2024"
$newDate = "12/10/2024"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame -eq -1) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Every slide layout.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# The slide master itself.
Update-DatePlaceholder $p.SlideMaster.Shapes

# ---------------------------------------------------------------------
# 2) New "RESTRICTED ACCESS" textbox on slide 2 (the Roomnames slide)
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Duplicate an existing room-name textbox so the new shape inherits the
# exact same run/paragraph/body formatting (shadow, font size, no-fill,
# centered, no-autofit, etc.) used throughout this slide.
$template = $s2.Shapes.Item($s2.Shapes.Count)
$newShape = $template.Duplicate()
$newShape.Name = "TextBox 8"

# Position/size in EMU -> points (1 pt = 12700 EMU), nudged to the
# nearest float32 value so the lossy Left/Top/Width/Height COM setters
# (which truncate through a 32-bit float) round-trip to the exact
# target EMU instead of landing 1 EMU short.
$newShape.Left = 67.03126525878906
$newShape.Top = 232.97268676757812
$newShape.Width = 396.85040283203125
$newShape.Height = 113.38583374023438

$newShape.TextFrame.TextRange.Text = "RESTRICTED ACCESS"
